$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the data/content that falls outside the new A1:C4 range,
# without disturbing the formatting (e.g. bold header style) of A1:C4.
$ws.Range("D1:E11").Clear()
$ws.Range("A5:C11").ClearContents()

# Write the new header row (keeps existing cell style/format)
$ws.Range("A1").Value = "Name"
$ws.Range("B1").Value = "Age"
$ws.Range("C1").Value = "City"

# Write the new data rows
$ws.Range("A2").Value = "Alice"
$ws.Range("B2").Value = 24
$ws.Range("C2").Value = "New York"

$ws.Range("A3").Value = "Bob"
$ws.Range("B3").Value = 30
$ws.Range("C3").Value = "Los Angeles"

$ws.Range("A4").Value = "Charlie"
$ws.Range("B4").Value = 22
$ws.Range("C4").Value = "Chicago"
